$wb = $excel.ActiveWorkbook

function Add-StockRow {
    param($ws, $row, $vals)

    $cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I")

    for ($i = 0; $i -lt 9; $i++) {
        $addr = $cols[$i] + $row
        if ($i -eq 7 -or $i -eq 8) {
            # Columns H (amount) and I (date) look numeric / date-like;
            # force text formatting so Excel keeps them as plain text,
            # matching the rest of the sheet, instead of silently
            # converting them to a number or a date serial value.
            $ws.Range($addr).NumberFormat = "@"
        }
        $ws.Range($addr).Value = $vals[$i]
    }
}

# Sheet "CY8C6247FTI-D52T" (worksheet #15): add row 2
$ws15 = $wb.Worksheets.Item(15)
$row = @("深圳市远洋乾坤电子科技有限公司", "notICCP", "SSCP", "CY8C6247FTI-D52T", "notSpotRanking", "notHotSell", "CYPRESS/赛普拉斯", "10000", "2023-02-20")
Add-StockRow $ws15 2 $row

# Sheet "CY8C6245LQI-S3D72" (worksheet #2): add row 2
$ws2 = $wb.Worksheets.Item(2)
$row = @("深圳市美信美科技有限公司", "notICCP", "SSCP", "CY8C6245LQI-S3D72", "notSpotRanking", "notHotSell", "Cypress", "5200", "2023-02-20")
Add-StockRow $ws2 2 $row

# Sheet "CY8C6247BZI-D34" (worksheet #5): add rows 3 and 4
$ws5 = $wb.Worksheets.Item(5)
$row = @("深圳市昇源芯科技有限公司", "notICCP", "SSCP", "CY8C6247BZI-D34", "notSpotRanking", "notHotSell", "CYPRESS/实单来谈特价出售", "3527", "2023-02-20")
Add-StockRow $ws5 3 $row
$row = @("深圳市振东芯电子科技有限公司", "notICCP", "SSCP", "CY8C6247BZI-D34", "notSpotRanking", "notHotSell", "CYPRESS/实单来谈特价出售", "3527", "2023-02-20")
Add-StockRow $ws5 4 $row

# Sheet "CY8C6247BZI-D44" (worksheet #6): add rows 3 and 4
$ws6 = $wb.Worksheets.Item(6)
$row = @("深圳市景宏锐科技有限公司", "notICCP", "notSSCP", "CY8C6247BZI-D44", "SpotRanking", "notHotSell", "CYPRESS", "5630", "2023-02-20")
Add-StockRow $ws6 3 $row
$row = @("深圳市联煌电子有限公司", "notICCP", "SSCP", "CY8C6247BZI-D44", "notSpotRanking", "notHotSell", "CYPRESS", "3790", "2023-02-20")
Add-StockRow $ws6 4 $row
